# REQ-27, REQ-28 and REQ-29 are added to system test report
#
# Adds three new test-case rows (28, 29, 30 -> spreadsheet rows 30-32) to the
# "Test Cases & Results" sheet, covering REQ-27, REQ-28 and REQ-29, then
# extends the conditional formatting / data validation / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")

# --- 1. Clone the formatting of the last populated row (29) onto the three
#        new rows so borders/fills/wrap-text/alignment match the existing
#        table rows exactly. ---
$ws.Range("B29:K29").Copy() | Out-Null
$ws.Range("B30:K30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B29:K29").Copy() | Out-Null
$ws.Range("B31:K31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B29:K29").Copy() | Out-Null
$ws.Range("B32:K32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Rows grew to the same auto-height (3 wrapped lines) as the other
# requirement rows.
$ws.Rows.Item(30).RowHeight = 43.2
$ws.Rows.Item(31).RowHeight = 43.2
$ws.Rows.Item(32).RowHeight = 43.2

# --- 2. Row 30 / TestCase 28 / REQ-27 ---
$ws.Range("B30").Formula = "=B29+1"
$ws.Range("D30").Value = "REQ-27"
$ws.Range("E30").Value = "High Impact"
$ws.Range("F30").Value = "Test if the loan status of the account in firebase is removed"
$ws.Range("G30").Value = "The return date must not be considered late from REQ-26"
$ws.Range("I30").Value = "Loan status of the test account in database is removed"
$ws.Range("J30").Value = "Loan status of the test account in database is removed"
$ws.Range("K30").Value = "Not Tested"

# --- 3. Row 31 / TestCase 29 / REQ-28 ---
$ws.Range("B31").Formula = "=B30+1"
$ws.Range("D31").Value = "REQ-28"
$ws.Range("E31").Value = "Low Impact"
$ws.Range("F31").Value = 'Test that LCD shows the confirmation messeage after the loan status is removed '
$ws.Range("G31").Value = "The user must go through REQ-27"
$ws.Range("I31").Value = "LCD shows the confirmation message"
$ws.Range("J31").Value = "LCD shows the confirmation message"
$ws.Range("K31").Value = "Not Tested"

# --- 4. Row 32 / TestCase 30 / REQ-29 ---
$ws.Range("B32").Formula = "=B31+1"
$ws.Range("D32").Value = "REQ-29"
$ws.Range("E32").Value = "High Impact"
$ws.Range("F32").Value = 'Test that the LCD shows "Please Scan Your Card" after the REQ-28 is completed '
$ws.Range("G32").Value = "The firebase is updated"
$ws.Range("H32").Value = "Follow the same steps as test case 29, no additional steps required "
$ws.Range("I32").Value = 'LCD displays "Please Scan Your Card"'
$ws.Range("J32").Value = 'LCD displays "Please Scan Your Card"'
$ws.Range("K32").Value = "Not Tested"

# --- 5. Extend the "Not Tested"/"Fail" conditional formatting and the
#        Pass/Fail/Not Tested data-validation dropdown down to row 32. ---
$ws.Range("K3:K29").FormatConditions.Delete()
$cf1 = $ws.Range("K3:K32").FormatConditions.Add(8, 3, '"Not Tested"')
$cf2 = $ws.Range("K3:K32").FormatConditions.Add(8, 3, '"Fail"')

$ws.Range("K3:K29").Validation.Delete()
$ws.Range("K3:K32").Validation.Add(3, 1, 1, "=Enums!$B$2:$B$4")

# --- 6. Scroll/selection state recorded by Excel after the edit. ---
$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H35").Select()
